$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-14 down to 10-15
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly price record
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Terminal La Palmera de La Serena"
$ws.Range("C9").Value = "Coquimbo"
$ws.Range("D9").Value = 44978
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = 100112039
$ws.Range("G9").Value = "Ciboulette"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 1000
$ws.Range("K9").Value = 1800
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 1900
$ws.Range("N9").Value = "`$/docena de atados"
$ws.Range("O9").Value = "Provincia del Elquí"
$ws.Range("P9").Value = 633
$ws.Range("Q9").Value = 3
$ws.Range("R9").Value = "Hortaliza"
